# Correção nos dados e inicio da analise PNAD 2009
#
# The sheet originally had two "section header" rows (row 5 and row 8)
# that only carried a row label ("situação do domicílio" and
# "grandes regiões e unidades da federação") with no numeric data next
# to them. Those rows are removed so the data rows flow contiguously,
# and the sub-header in B2 (which incorrectly held the literal text
# "unnamed: 1_level_1") is corrected to read "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mislabeled sub-header cell.
$ws.Range("B2").Value = "total"

# Remove the empty "situação do domicílio" and "grandes regiões e
# unidades da federação" section rows, working bottom-up so row
# numbers of not-yet-deleted rows stay valid.
$ws.Rows("8").Delete()
$ws.Rows("5").Delete()

Write-Host "applied PNAD 2009 roubo cv124021a corrections"
